$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CapitalCall")

# Add two new "send notification" flag columns after the existing
# "Unit Price/Premium" column (K): L = Send Payment Notification,
# M = Send Call Notice.
$ws.Range("L1").Value = "Send Payment Notification"
$ws.Range("M1").Value = "Send Call Notice"

# Every existing capital-call row defaults both new flags to "Yes".
$ws.Range("L2:L4").Value = "Yes"
$ws.Range("M2:M4").Value = "Yes"

# Reflect the updated view/selection now that the used range grew to
# include the new columns.
$ws.Range("K1").Select()
